$p = $ppt.ActivePresentation

# --- Slide 10 (sldId 276): "Title 1" shape (id=6) ---
# Move/resize the title textbox and update its title text to mention revenue.
$s10 = $p.Slides.Item(10)
$t10 = $s10.Shapes.Item(2)

$t10.Left = 306.810787
$t10.Top = -5.837953
$t10.Width = 364.540551
$t10.Height = 41.074341

$t10.TextFrame.TextRange.Text = "Genre vs rating and revenue"

# --- Slide 11 (sldId 277): "Title 1" shape (id=6) ---
# Move/resize the title textbox and rename it to the ROI slide title.
$s11 = $p.Slides.Item(11)
$t11 = $s11.Shapes.Item(2)

$t11.Left = 280.125847
$t11.Top = -5.837884
$t11.Width = 423.009301
$t11.Height = 41.074341

$t11.TextFrame.TextRange.Text = "Genre vs Return on investment"
